# The two species records in rows 3 and 4 of the "Artfynd" sheet were
# swapped (row 3 now contains what used to be row 4's data, and vice
# versa). Columns whose values are identical between the two rows
# (D, T, U, V, W, Y, AA, AD, AE, AG, AT, AY, I) are left untouched.
# Row 3 additionally loses its Z/AB ("time") values and gains row 4's
# AC ("comment") value, while row 4 gains the Z/AB values and loses AC.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# New row 3 values (previously on row 4)
$ws.Range("A3").Value = 130861152
$ws.Range("B3").Value = 91804
$ws.Range("E3").Value = 1108
$ws.Range("F3").Value = "Harticka"
$ws.Range("G3").Value = "Pelloporus leporinus"
$ws.Range("H3").Value = "(Fr.) Krieglst."
$ws.Range("P3").Value = "Djupbäcken, Jmt"
$ws.Range("Q3").Value = 442868
$ws.Range("R3").Value = 7039767
$ws.Range("S3").Value = 10
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").Value = "I stående levande gran med full längd."
$ws.Range("AW3").Value = "Kristian Zackrisson"
$ws.Range("AX3").Value = "Kristian Zackrisson"

# New row 4 values (previously on row 3)
$ws.Range("A4").Value = 130853761
$ws.Range("B4").Value = 79244
$ws.Range("E4").Value = 230405
$ws.Range("F4").Value = "Garnlav (ssp. sarmentosa)"
$ws.Range("G4").Value = "Alectoria sarmentosa subsp. sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("P4").Value = "Djupbäcken, Djupbäcken, Jmt"
$ws.Range("Q4").Value = 442771
$ws.Range("R4").Value = 7039709
$ws.Range("S4").Value = 20
$ws.Range("Z4").Value = "11:05"
$ws.Range("AB4").Value = "11:05"
$ws.Range("AC4").ClearContents()
$ws.Range("AW4").Value = "Maria Danvind"
$ws.Range("AX4").Value = "Maria Danvind"
